$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.971.95"
$ws.Range("E2").Value = "  +1.67%  "

$ws.Range("D3").Value = "2.244.10"
$ws.Range("E3").Value = "  +0.74%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.09"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.27"
$ws.Range("E6").Value = "  +1.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.574"
$ws.Range("E7").Value = "  -1.20%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.544"
$ws.Range("E9").Value = "  -2.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.70"
$ws.Range("E10").Value = "  -0.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0825"
$ws.Range("E11").Value = "  -0.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.50"
$ws.Range("E12").Value = "  -2.54%  "

$ws.Range("E13").Value = "  -1.66%  "

$ws.Range("D14").Value = "2.586.33"
$ws.Range("E14").Value = "  +0.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.848"
$ws.Range("E15").Value = "  -1.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.26"
$ws.Range("E16").Value = "  -0.51%  "

$ws.Range("D17").Value = "2.244.73"
$ws.Range("E17").Value = "  +0.59%  "

$ws.Range("D18").Value = "43.895.81"
$ws.Range("E18").Value = "  +1.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.30"
$ws.Range("E19").Value = "  -4.91%  "

$ws.Range("D20").Value = "0.0₃0974"
$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.41"
$ws.Range("E21").Value = "  -2.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.19"
$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("E23").Value = "  -3.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "233.31"
$ws.Range("E24").Value = "  -1.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.06"
$ws.Range("E25").Value = "  -5.50%  "

$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.60"
$ws.Range("E27").Value = "  +5.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.76"
$ws.Range("E28").Value = "  +5.60%  "

$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.03"
$ws.Range("E30").Value = "  -5.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.58"
$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.04"
$ws.Range("E32").Value = "  -0.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0840"
$ws.Range("E33").Value = "  -2.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.68"
$ws.Range("E34").Value = "  +0.28%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.11"
$ws.Range("E35").Value = "  -4.75%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  +7.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.92"
$ws.Range("E37").Value = "  +4.90%  "

$ws.Range("E38").Value = "  -1.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.22"
$ws.Range("E39").Value = "  +12.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.66"
$ws.Range("E40").Value = "  -1.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.15"
$ws.Range("E41").Value = "  -5.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0313"
$ws.Range("E42").Value = "  -1.76%  "

$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").Value = "1.761.30"
$ws.Range("E44").Value = "  -0.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "74.41"
$ws.Range("E45").Value = "  +1.25%  "

$ws.Range("E46").Value = "  -3.77%  "

$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "80.88"
$ws.Range("E47").Value = "  -3.20%  "

$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.16"
$ws.Range("E48").Value = "  -2.25%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.00"
$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.66"
$ws.Range("E50").Value = "  +1.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "57.13"
$ws.Range("E51").Value = "  -1.52%  "
